$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename existing "Passives" column (S) to "PassiveModifier" semantics
$ws.Range("S1").Value = "PassiveModifier"
$ws.Range("S3").Value = "被动数值"
$ws.Range("S4").Value = 1

# Add two new trailing columns: Passive (X) and PassiveDesc (Y)
$ws.Range("X1").Value = "Passive"
$ws.Range("X2").Value = "string"
$ws.Range("X3").Value = "被动"
$ws.Range("X4").Value = "避柳"

$ws.Range("Y1").Value = "PassiveDesc"
$ws.Range("Y2").Value = "string"
$ws.Range("Y3").Value = "被动描述"
$ws.Range("Y4").Value = "场上有单位进入危险区时减少自身*PassiveModifier*费用。一次行动只能触发一次。（直到下个回合无法触发）"

# Update selection / view state to match target
$ws.Range("Y4").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
